$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AVA")

# Append the two new processes to the exclusion list held in column T
# (process) for every grid-node row (T11:T159). All of these cells
# originally shared the identical text "e_demand,ev_battery".
$ws.Range("T11:T159").Value = "e_demand,ev_battery,H2prd_Elc_PEM,H2prd_Elc_ALK"

# The longer text no longer fits the previous best-fit width, so widen
# column T (20) to match the width already used by the analogous
# "process" columns E (5) and J (10), dropping the best-fit auto sizing.
$ws.Columns.Item(20).ColumnWidth = $ws.Columns.Item(5).ColumnWidth
